$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lines")

# Update data values: D2/D3 125 -> 0.5, F2 0.083 -> 0.0083
$ws.Range("D2").Value = 0.5
$ws.Range("F2").Value = 0.0083
$ws.Range("D3").Value = 0.5

# Move active selection on the sheet to F3
$ws.Activate()
$ws.Range("F3").Select()
